$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.930.34'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '3.393.60'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '572.19'
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.82'
$ws.Range('E6').Value = '  -1.84%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.394.18'
$ws.Range('E7').Value = '  -1.35%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.55'
$ws.Range('E10').Value = '  -1.66%  '
$ws.Range('E11').Value = '  -2.11%  '
$ws.Range('E12').Value = '  +2.06%  '
$ws.Range('D13').Value = '3.971.75'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('E14').Value = '  +2.08%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '28.20'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('D17').Value = '3.386.25'
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').Value = '60.994.39'
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.87'
$ws.Range('E20').Value = '  -2.91%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '8.99'
$ws.Range('E21').Value = '  -4.55%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '383.40'
$ws.Range('E22').Value = '  -3.78%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.560'
$ws.Range('E23').Value = '  -1.52%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '74.37'
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000118'
$ws.Range('E26').Value = '  -4.23%  '
$ws.Range('D27').Value = '3.525.72'
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  -2.94%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.99'
$ws.Range('E31').Value = '  -3.20%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.15'
$ws.Range('E32').Value = '  -1.70%  '
$ws.Range('E33').Value = '  -2.55%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '23.57'
$ws.Range('E35').Value = '  -1.66%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.02'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '167.25'
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('D38').Value = '3.421.72'
$ws.Range('E38').Value = '  -1.28%  '
$ws.Range('E39').Value = '  -2.63%  '
$ws.Range('E40').Value = '  -4.10%  '
$ws.Range('E41').Value = '  -1.75%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '27.69'
$ws.Range('E42').Value = '  +2.53%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.781'
$ws.Range('E43').Value = '  -2.60%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '42.17'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.44'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('E47').Value = '  -3.35%  '
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('D49').Value = '2.481.30'
$ws.Range('E49').Value = '  -4.60%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.83'
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '22.98'
$ws.Range('E51').Value = '  -0.57%  '
